# Apply cryptocurrency price/volume updates (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.071.73"
$ws.Range("E2").Value = "  -3.02%  "

$ws.Range("D3").Value = "3.508.36"
$ws.Range("E3").Value = "  -4.85%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.29"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.11"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.41%  "

$ws.Range("D7").Value = "3.501.19"
$ws.Range("E7").Value = "  -4.84%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.606"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.06%  "

$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  -5.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.49"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("E12").Value = "  -4.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.56"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.77%  "

$ws.Range("E14").Value = "  -4.45%  "

$ws.Range("D15").Value = "4.078.61"
$ws.Range("E15").Value = "  -4.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.53"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "623.49"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -7.86%  "

$ws.Range("D18").Value = "69.115.32"
$ws.Range("E18").Value = "  -3.14%  "

$ws.Range("D19").Value = "3.509.05"
$ws.Range("E19").Value = "  -4.95%  "

$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.39"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.83%  "

$ws.Range("E22").Value = "  -3.80%  "

$ws.Range("E23").Value = "  -6.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.88"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -8.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.16"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.79"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.45%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("E28").Value = "  -6.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.31"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -8.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.54"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.99%  "

$ws.Range("E31").Value = "  -7.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.51"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.75%  "

$ws.Range("E33").Value = "  -8.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.97"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "633.65"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.72"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.06%  "

$ws.Range("E37").Value = "  -5.26%  "

$ws.Range("E38").Value = "  -15.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.58"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0446"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.137"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.56%  "

$ws.Range("D43").Value = "3.364.72"
$ws.Range("E43").Value = "  -8.43%  "

$ws.Range("E44").Value = "  -6.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.81"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.91%  "

$ws.Range("D46").Value = "0.0₃0688"
$ws.Range("E46").Value = "  -9.80%  "

$ws.Range("E47").Value = "  -7.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.76"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.61%  "

$ws.Range("E49").Value = "  -2.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.39"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.66"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +14.88%  "
